$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (workbook stores them as inline strings,
# not numbers) so values like '15.60' or '0.0000277' keep their exact formatting
# instead of being auto-coerced to numeric literals by COM's Value setter.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{Row=2; D='70.653.89'; E='  -0.38%  '},
    @{Row=3; D='3.518.59'; E='  -2.19%  '},
    @{Row=4; E='  -0.06%  '},
    @{Row=5; D='621.20'; E='  +3.34%  '},
    @{Row=6; D='172.61'; E='  -0.85%  '},
    @{Row=7; D='0.611'; E='  -1.18%  '},
    @{Row=8; D='3.514.77'; E='  -2.12%  '},
    @{Row=10; E='  -0.34%  '},
    @{Row=11; D='7.16'; E='  -3.86%  '},
    @{Row=12; E='  -0.80%  '},
    @{Row=13; D='46.36'; E='  -1.55%  '},
    @{Row=14; D='0.0000277'; E='  -1.05%  '},
    @{Row=15; D='4.096.95'; E='  -2.00%  '},
    @{Row=16; E='  -1.32%  '},
    @{Row=17; D='608.93'; E='  -1.20%  '},
    @{Row=18; D='3.526.78'; E='  -1.97%  '},
    @{Row=19; D='70.769.74'; E='  -0.31%  '},
    @{Row=20; E='  +1.24%  '},
    @{Row=21; D='17.79'; E='  +1.55%  '},
    @{Row=22; E='  -0.88%  '},
    @{Row=23; D='9.13'; E='  -1.44%  '},
    @{Row=24; D='15.60'; E='  -3.05%  '},
    @{Row=25; D='97.62'; E='  -0.05%  '},
    @{Row=26; D='3.73'; E='  -1.52%  '},
    @{Row=27; E='  +0.07%  '},
    @{Row=28; E='  -3.76%  '},
    @{Row=29; D='33.56'; E='  -1.62%  '},
    @{Row=30; D='9.07'; E='  -2.01%  '},
    @{Row=31; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='3.01'; E='  -2.38%  '},
    @{Row=32; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='8.07'; E='  -5.29%  '},
    @{Row=33; E='  -0.90%  '},
    @{Row=34; D='632.33'; E='  -1.55%  '},
    @{Row=35; D='6.83'; E='  -6.01%  '},
    @{Row=36; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='10.81'; E='  -0.73%  '},
    @{Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0991'; E='  -2.54%  '},
    @{Row=38; D='0.0484'; E='  +0.85%  '},
    @{Row=39; D='3.43'; E='  -8.39%  '},
    @{Row=40; D='56.66'; E='  -1.35%  '},
    @{Row=41; E='  +0.23%  '},
    @{Row=42; D='0.142'; E='  -0.37%  '},
    @{Row=43; D='3.348.81'; E='  -2.06%  '},
    @{Row=44; D='0.0₃0721'; E='  +0.17%  '},
    @{Row=45; E='  -0.56%  '},
    @{Row=46; E='  -4.33%  '},
    @{Row=47; D='31.88'; E='  -3.88%  '},
    @{Row=48; E='  -6.01%  '},
    @{Row=49; E='  -0.73%  '},
    @{Row=50; D='133.98'; E='  +0.75%  '},
    @{Row=51; E='  +0.02%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey('B')) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
